# Update database: rial_cumulative.xlsx (darou/desobhan income quarterly)
# - Refresh the "12 ماهه منتهی به 1401/12" (12-month cumulative) period's
#   publish-date labels and financial figures in column M (and the matching
#   header label in column I) to reflect the new report date 1402-03-07.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Publish-date header labels (row 9)
$ws.Range("I9").Value = "1402-03-07 (8)"
$ws.Range("M9").Value = "1402-03-07 (2)"

# Updated cumulative figures for column M (12 ماهه منتهی به 1401/12)
$ws.Range("M12").Value = -6279464
$ws.Range("M13").Value = 6075983
$ws.Range("M14").Value = -654396
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 75335
$ws.Range("M17").Value = 5496922
$ws.Range("M19").Value = 91530
$ws.Range("M20").Value = 4576816
$ws.Range("M21").Value = -750790
$ws.Range("M22").Value = 3826026
$ws.Range("M24").Value = 3826026
$ws.Range("M25").Value = 1054
$ws.Range("M27").Value = 1054

$wb.Save()
